# "corrected story points from 18 to 21"
# Rows 21-23 (Reference Story No. 19, 20, 21) had a Pool/story-points value
# of 18 that should actually be 21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = 21
$ws.Range("C22").Value = 21
$ws.Range("C23").Value = 21

# Leave the view scrolled down to / focused on the corrected rows, matching
# where the author was working when they made the fix.
$ws.Range("C21").Select()
